# Append the new September schedule entries (rows 12-20) to Sheet1,
# reusing the existing "INITIALS/DATE/COURSE_ID/TIME(MINUTES)" columns,
# then leave the selection on D19 (last entered value) as in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "m/d/yy"   # matches the existing short-date format used in column B

$rows = @(
    @{ Row = 12; Initials = "JG"; Date = 45181; Course = "DS160-01";   Minutes = 75 },
    @{ Row = 13; Initials = "JG"; Date = 45181; Course = "MATH205-03"; Minutes = 90 },
    @{ Row = 14; Initials = "JG"; Date = 45182; Course = "ECON110-03"; Minutes = 75 },
    @{ Row = 15; Initials = "JG"; Date = 45183; Course = "DS160-01";   Minutes = 90 },
    @{ Row = 16; Initials = "JG"; Date = 45183; Course = "MATH205-03"; Minutes = 90 },
    @{ Row = 17; Initials = "JG"; Date = 45187; Course = "ECON110-03"; Minutes = 90 },
    @{ Row = 18; Initials = "JG"; Date = 45187; Course = "THEO200-05"; Minutes = 90 },
    @{ Row = 19; Initials = "JG"; Date = 45188; Course = "DS160-01";   Minutes = 75 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Initials
    $cellB = $ws.Cells.Item($r.Row, 2)
    $cellB.NumberFormat = $dateFmt
    $cellB.Value = $r.Date
    $ws.Cells.Item($r.Row, 3).Value = $r.Course
    $ws.Cells.Item($r.Row, 4).Value = $r.Minutes
}

# Row 20 only has the initials filled in.
$ws.Cells.Item(20, 1).Value = "JG"

$ws.Range("D19").Select()
